$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the "sub.google.com" header row and "x" value row from column BS (71) to CC (81)
for ($c = 72; $c -le 81; $c++) {
    $ws.Cells.Item(1, $c).Value = "sub.google.com"
    $ws.Cells.Item(2, $c).Value = "x"
}

# Update the frozen pane top-left cell and the active selection in the bottom-right pane
$ws.Activate()
$ws.Range("BS2").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("BY3").Select()
